# Revise config file handling
# Append a new row (row 47) of device/config data to each of the four
# worksheets (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2).

$wb = $excel.ActiveWorkbook

function Add-ConfigRow {
    param(
        $ws,
        [double]$timeVal,
        [string]$col_b,
        [string]$col_c,
        [string]$col_d,
        [string]$col_e,
        [double]$col_f,
        [string]$col_g,
        [double]$col_h,
        [double]$col_i
    )

    $row = 47

    $ws.Cells.Item($row, 1).Value = $timeVal
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $col_b
    $ws.Cells.Item($row, 3).Value = $col_c
    $ws.Cells.Item($row, 4).Value = $col_d
    $ws.Cells.Item($row, 5).Value = $col_e

    $ws.Cells.Item($row, 6).Value = $col_f
    $ws.Cells.Item($row, 7).Value = [double]$col_g
    $ws.Cells.Item($row, 8).Value = $col_h
    $ws.Cells.Item($row, 9).Value = $col_i
}

# Sheet 1: MID_LFT_#1
$ws1 = $wb.Worksheets.Item("MID_LFT_#1")
Add-ConfigRow $ws1 45833.46165509259 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x70" "0x07" 400 "5.68631262647113e+23" 368 7

# Sheet 2: MID_LFT_#2
$ws2 = $wb.Worksheets.Item("MID_LFT_#2")
Add-ConfigRow $ws2 45833.46165509259 "0x01,0x7c" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x64" "0x19" 380 "5.68432987514711e+23" 356 25

# Sheet 3: MID_PLT_#1
$ws3 = $wb.Worksheets.Item("MID_PLT_#1")
Add-ConfigRow $ws3 45833.46165509259 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x00,0x6A" "0x15" 110 "5.68631262647113e+23" 106 15

# Sheet 4: MID_PLT_#2
$ws4 = $wb.Worksheets.Item("MID_PLT_#2")
Add-ConfigRow $ws4 45833.46165509259 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x00,0x7F" "0x9" 130 "5.68631262647113e+23" 127 9

Write-Output "Added row 47 to all 4 sheets"
